$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B29").Value = "Method Long aroundLongMethod(String str)"
$ws.Range("B30").Value = "return Long.valueOf(str);"
$ws.Range("B33").Value = "Method void aroundVoidMethod()"
$ws.Range("B34").Value = "return;"

$ws.Range("B33").Select()
